# Scheduled runner update: refresh market-board derived figures across the
# per-job Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1963.1666
$ws.Range("J17").Value = 2055.8
$ws.Range("L17").Value = 6167.400000000001
$ws.Range("N17").Value = -6503.400000000001

$ws.Range("H41").Value = 1808.3077
$ws.Range("I41").Value = 208.2
$ws.Range("K41").Value = 208.2
$ws.Range("M41").Value = 231.8

$ws.Range("H76").Value = 3711.8572
$ws.Range("I76").Value = 3750
$ws.Range("K76").Value = 3750
$ws.Range("M76").Value = -3435

$ws.Range("H79").Value = 3711.8572
$ws.Range("I79").Value = 3750
$ws.Range("K79").Value = 3750
$ws.Range("M79").Value = -2658

$ws.Range("H86").Value = 127075
$ws.Range("I86").Value = 168766.67
$ws.Range("K86").Value = 168766.67
$ws.Range("M86").Value = -167643.67

$ws.Range("H89").Value = 127075
$ws.Range("I89").Value = 168766.67
$ws.Range("K89").Value = 843833.3500000001
$ws.Range("M89").Value = -838217.3500000001

$ws.Range("H106").Value = 2359.375
$ws.Range("J106").Value = 2000
$ws.Range("L106").Value = 2000
$ws.Range("N106").Value = -3262

$ws.Range("H132").Value = 5874.1377
$ws.Range("I132").Value = 1838
$ws.Range("K132").Value = 5514
$ws.Range("M132").Value = -2984

# --- ARM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20093.543
$ws.Range("I32").Value = 19300.219
$ws.Range("J32").Value = 26598.8
$ws.Range("K32").Value = 19300.219
$ws.Range("L32").Value = 26598.8
$ws.Range("M32").Value = -19013.219
$ws.Range("N32").Value = -27172.8

$ws.Range("H43").Value = 19997.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 19997.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 19997.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -20623.5

$ws.Range("H92").Value = 4083920
$ws.Range("J92").Value = 4083920
$ws.Range("L92").Value = 4083920
$ws.Range("N92").Value = -4088912

$ws.Range("H109").Value = 51158.285
$ws.Range("J109").Value = 51158.285
$ws.Range("L109").Value = 51158.285
$ws.Range("N109").Value = -53932.285

# --- BSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5437.095
$ws.Range("I20").Value = 5583.4287
$ws.Range("J20").Value = 5144.4287
$ws.Range("K20").Value = 5583.4287
$ws.Range("L20").Value = 5144.4287
$ws.Range("M20").Value = -5336.4287
$ws.Range("N20").Value = -5638.4287

$ws.Range("H86").Value = 3871.7058
$ws.Range("I86").Value = 2632.1
$ws.Range("J86").Value = 5642.5713
$ws.Range("K86").Value = 2632.1
$ws.Range("L86").Value = 5642.5713
$ws.Range("M86").Value = -1509.1
$ws.Range("N86").Value = -7888.5713

$ws.Range("H89").Value = 3871.7058
$ws.Range("I89").Value = 2632.1
$ws.Range("J89").Value = 5642.5713
$ws.Range("K89").Value = 13160.5
$ws.Range("L89").Value = 28212.8565
$ws.Range("M89").Value = -7544.5
$ws.Range("N89").Value = -39444.85649999999

$ws.Range("H134").Value = 2443.6333
$ws.Range("I134").Value = 2009.28
$ws.Range("J134").Value = 4615.4
$ws.Range("K134").Value = 6027.84
$ws.Range("L134").Value = 13846.2
$ws.Range("M134").Value = -3492.84
$ws.Range("N134").Value = -18916.2

# --- CRP ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2593.3699
$ws.Range("I31").Value = 2256.4854
$ws.Range("K31").Value = 2256.4854
$ws.Range("M31").Value = -1961.4854

$ws.Range("H34").Value = 2593.3699
$ws.Range("I34").Value = 2256.4854
$ws.Range("K34").Value = 2256.4854
$ws.Range("M34").Value = -2054.4854

$ws.Range("H134").Value = 2470.8286
$ws.Range("I134").Value = 1545.4584
$ws.Range("K134").Value = 4636.3752
$ws.Range("M134").Value = -2101.3752

$ws.Range("H141").Value = 119124.66
$ws.Range("J141").Value = 126407.52
$ws.Range("L141").Value = 126407.52
$ws.Range("N141").Value = -136767.52

# --- CUL ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 791.7692
$ws.Range("I5").Value = 304.16666
$ws.Range("J5").Value = 1209.7142
$ws.Range("K5").Value = 912.4999799999999
$ws.Range("L5").Value = 3629.1426
$ws.Range("M5").Value = -800.4999799999999
$ws.Range("N5").Value = -3853.1426

$ws.Range("H80").Value = 5749.75
$ws.Range("J80").Value = 5749.75
$ws.Range("L80").Value = 17249.25
$ws.Range("N80").Value = -19121.25

$ws.Range("H83").Value = 5749.75
$ws.Range("J83").Value = 5749.75
$ws.Range("L83").Value = 51747.75
$ws.Range("N83").Value = -61107.75

$ws.Range("H135").Value = 791.7692
$ws.Range("I135").Value = 304.16666
$ws.Range("J135").Value = 1209.7142
$ws.Range("K135").Value = 2737.49994
$ws.Range("L135").Value = 10887.4278
$ws.Range("M135").Value = -202.4999399999997
$ws.Range("N135").Value = -15957.4278

# --- GSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 25000
$ws.Range("I26").Value = 25000
$ws.Range("K26").Value = 25000
$ws.Range("M26").Value = -24720

$ws.Range("H50").Value = 25000
$ws.Range("I50").Value = 25000
$ws.Range("K50").Value = 25000
$ws.Range("M50").Value = -24502

# --- LTW ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2701.4119
$ws.Range("I16").Value = 1648.0769
$ws.Range("K16").Value = 1648.0769
$ws.Range("M16").Value = -1478.0769

$ws.Range("H40").Value = 25645932
$ws.Range("I40").Value = 3973.375
$ws.Range("J40").Value = 66673064
$ws.Range("K40").Value = 3973.375
$ws.Range("L40").Value = 66673064
$ws.Range("M40").Value = -3837.375
$ws.Range("N40").Value = -66673336

$ws.Range("H46").Value = 3794.2273
$ws.Range("J46").Value = 3982.0557
$ws.Range("L46").Value = 3982.0557
$ws.Range("N46").Value = -4358.0557

$ws.Range("H114").Value = 52484.5
$ws.Range("J114").Value = 52484.5
$ws.Range("L114").Value = 52484.5
$ws.Range("N114").Value = -61162.5

$ws.Range("H122").Value = 6720.744
$ws.Range("I122").Value = 3897.4614
$ws.Range("J122").Value = 11038.706
$ws.Range("K122").Value = 11692.3842
$ws.Range("L122").Value = 33116.118
$ws.Range("M122").Value = -9242.3842
$ws.Range("N122").Value = -38016.118

# --- WVR ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H81").Value = 12589.789
$ws.Range("I81").Value = 5527.857
$ws.Range("J81").Value = 16709.25
$ws.Range("K81").Value = 11055.714
$ws.Range("L81").Value = 33418.5
$ws.Range("M81").Value = -9994.714
$ws.Range("N81").Value = -35540.5

$ws.Range("H84").Value = 12589.789
$ws.Range("I84").Value = 5527.857
$ws.Range("J84").Value = 16709.25
$ws.Range("K84").Value = 55278.57
$ws.Range("L84").Value = 167092.5
$ws.Range("M84").Value = -49974.57
$ws.Range("N84").Value = -177700.5

$ws.Range("H132").Value = 1620.2716
$ws.Range("I132").Value = 768.73914
$ws.Range("J132").Value = 1957.9482
$ws.Range("K132").Value = 2306.21742
$ws.Range("L132").Value = 5873.8446
$ws.Range("M132").Value = 223.7825800000001
$ws.Range("N132").Value = -10933.8446

$ws.Range("H138").Value = 98049.5
$ws.Range("I138").Value = 98000
$ws.Range("K138").Value = 98000
$ws.Range("M138").Value = -92860

$ws.Range("H140").Value = 92909.2
$ws.Range("J140").Value = 92909.2
$ws.Range("L140").Value = 92909.2
$ws.Range("N140").Value = -103269.2
